# Generate Report for Handback
# Updates the "zh-cn" and "de-de" worksheets with the latest handback
# information for the 3200f3a5-... file (row 5) and refreshes a couple of
# downstream xliff-generation timestamps (row 6), mirroring a "handback
# report" re-generation.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d46c77ae4e86078b016c8510e06e9eae8eb01cd5/e2e/3200f3a5-b483-4f30-a445-552926e36023.md"
$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d7d20b8d429445fb88fa003928b42041a6c7e77/e2e/3200f3a5-b483-4f30-a445-552926e36023.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d46c77ae4e86078b016c8510e06e9eae8eb01cd5/e2e/3200f3a5-b483-4f30-a445-552926e36023.md."

function Update-LocalizationSheet {
    param(
        [string]$SheetName,
        [string]$TargetXlfName,
        [string]$HandbackDateTime,
        [string]$G6Value,
        [string]$H6Value
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen columns I (Latest Target File), J (Latest Handback File) and
    # P (Error Detail) to 40 characters, same as the other "wide" columns.
    $ws.Range("I1").EntireColumn.ColumnWidth = 39.17
    $ws.Range("J1").EntireColumn.ColumnWidth = 39.17
    $ws.Range("P1").EntireColumn.ColumnWidth = 39.17

    # Remember the current hyperlinks (in sheet order) so we can rebuild
    # them in the correct order after inserting the new one.
    $hlRefs = @()
    $hlAddrs = @()
    $hlDisplays = @()
    foreach ($hl in $ws.Hyperlinks) {
        $hlRefs += $hl.Range.Address()
        $hlAddrs += $hl.Address
        $hlDisplays += $hl.TextToDisplay
    }

    # Row 5: the handback for 3200f3a5-... has now been received - record
    # the latest target/handback file names and datetime, and flag that the
    # handback version is not the latest.
    $ws.Range("J5").Value = $TargetXlfName
    $ws.Range("K5").Value = $HandbackDateTime
    $ws.Range("P5").Value = $errorMsg

    # Row 6: refresh the generated xliff file names for the next file.
    $ws.Range("G6").Value = $G6Value
    $ws.Range("H6").Value = $H6Value

    # Rebuild hyperlinks: keep A2..A5 as-is, insert a new one on I5 (points
    # to the same handback markdown file as A5), then re-add the A6 link.
    $ws.Range("A1").Hyperlinks.Delete()

    for ($i = 0; $i -lt $hlRefs.Count - 1; $i++) {
        $ws.Hyperlinks.Add($ws.Range($hlRefs[$i]), $hlAddrs[$i], "", "", $hlDisplays[$i]) | Out-Null
    }

    $ws.Hyperlinks.Add($ws.Range("I5"), $targetUrl, "", "", "3200f3a5-b483-4f30-a445-552926e36023.md") | Out-Null
    $ws.Range("I5").Font.Underline = $true
    $ws.Range("I5").Font.Color = 6710373

    $lastIdx = $hlRefs.Count - 1
    $ws.Hyperlinks.Add($ws.Range($hlRefs[$lastIdx]), $hlAddrs[$lastIdx], "", "", $hlDisplays[$lastIdx]) | Out-Null
}

Update-LocalizationSheet -SheetName "zh-cn" `
    -TargetXlfName "3200f3a5-b483-4f30-a445-552926e36023.159151e919cf2d883bcb39c31250c5a48beed386.zh-cn.xlf" `
    -HandbackDateTime "2016-10-18 02:52:09" `
    -G6Value "ead9b05c-6df6-4f2d-9561-2cf7d1a36e36.5d4e3e047f66d9fd21175eb097ca7560e4ec03c0.zh-cn.xlf" `
    -H6Value "2016-10-18 02:48:00"

Update-LocalizationSheet -SheetName "de-de" `
    -TargetXlfName "3200f3a5-b483-4f30-a445-552926e36023.159151e919cf2d883bcb39c31250c5a48beed386.de-de.xlf" `
    -HandbackDateTime "2016-10-18 02:52:47" `
    -G6Value "ead9b05c-6df6-4f2d-9561-2cf7d1a36e36.5d4e3e047f66d9fd21175eb097ca7560e4ec03c0.de-de.xlf" `
    -H6Value "2016-10-18 02:48:23"

Write-Output "Done updating handback report."
